$d = $word.ActiveDocument

# The "Place approximately <N> feet of  FIBER cable." sentence currently reads
# "13664" (rendered from two adjacent runs containing "136" and "64").
# Update the footage figure to "13654" (runs containing "13" and "654").
# Find/Execute operates over the document's logical text stream, so it matches
# across the run boundary even though the digits live in separate <w:r> elements.
$found = $d.Content.Find.Execute("13664", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "13654", 2)

if (-not $found) {
    throw "Could not find the '13664' footage figure to replace."
}
